$d = $word.ActiveDocument

# Locate the paragraph that ends the document body: "Signatur des Absenders ..."
$rng = $d.Content
$found = $rng.Find.Execute("Signatur des Absenders bestätigt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $rng.Paragraphs(1)

# 1) New plain paragraph right after it
$targetPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs($d.Paragraphs.Count)
$p1.Range.Text = "Mail selbst wurde mit dem PublicKey des Empfängers verschlüsselt und kann nur mit dessen privateKey entschlüsselt werden"

# 2) Empty paragraph
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($d.Paragraphs.Count)

# 3) Bold + underlined heading-like paragraph (built from several sentence pieces)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($d.Paragraphs.Count)

# 4) Trailing empty paragraph (inserted now, before bold/underline formatting
#    is applied to p3, so it stays unformatted)
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs($d.Paragraphs.Count)

$p3.Range.Text = "Beim Senden einer Mail muss der publicKey des Empfängers bei der Uni angefragt werden!"
$p3.Range.Font.Bold = 1
$p3.Range.Font.BoldBi = 1
$p3.Range.Font.Underline = 1

Write-Host "Paragraphs now: $($d.Paragraphs.Count)"
